$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Agropecuária - Bambuí): Inscritos 25->27, Pagos 14->15, Inscrições homologadas 16->17
$ws.Range("E2").Value = 27
$ws.Range("F2").Value = 15
$ws.Range("H2").Value = 17

# Row 7 (Mineração - Congonhas): Inscritos 31->32
$ws.Range("E7").Value = 32

# Row 16 (Logística - Ribeirão das Neves): Inscritos 313->318
$ws.Range("E16").Value = 318

# Row 18 (Segurança do Trabalho - Santa Luzia): Inscritos 96->97
$ws.Range("E18").Value = 97
